$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) The "Run" column (E2:E17) all flip from TRUE to FALSE.
$ws.Range("E2:E17").Value = $false

# 2) Data validation on column E currently covers the range via three
#    separate pieces (E8, E2:E7, E9:E1048576) - collapse it into one
#    contiguous rule E2:E1048576 with the same list source.
$ws.Range("E2:E1048576").Validation.Delete()
$ws.Range("E2:E1048576").Validation.Add(3, 1, 1, '"TRUE,FALSE"')

# 3) Move the active selection from E9 to E12.
$ws.Range("E12").Select()

"done"
